$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT, without Excel's
# automatic "looks like a date/number" re-typing kicking in (which would
# otherwise turn e.g. "2025-11-09" into a date serial number).
#
# Trick: put the text in via a formula that evaluates to a string
# ("2025-11-09" -> ="2025-11-09"), then Copy / PasteSpecial(values) the
# cell onto itself. Paste-as-values keeps the already-computed string but
# drops the formula, and (unlike typing text / quote-prefixing / setting
# NumberFormat="@") it does not allocate a new text-forced cell style, so
# the result matches a plain text cell with the default style.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Existing rows 2-5: Amount column was stored as text ("1000" etc.);
#     it should become a genuine number --------------------------------
$ws.Range("B2").Value = 1000
$ws.Range("B3").Value = 12000
$ws.Range("B4").Value = 2000
$ws.Range("B5").Value = 2000

# --- New rows 6-10 -------------------------------------------------------
$ws.Range("A6").Value = "shopping"
$ws.Range("B6").Value = 1000
Set-LiteralText $ws.Cells.Item(6, 3) "2025-11-09"

$ws.Range("A7").Value = "movies"
$ws.Range("B7").Value = 2000
Set-LiteralText $ws.Cells.Item(7, 3) "2025-11-09"

$ws.Range("A8").Value = "groceries"
$ws.Range("B8").Value = 1200
Set-LiteralText $ws.Cells.Item(8, 3) "2025-11-09"

$ws.Range("A9").Value = "Rent"
$ws.Range("B9").Value = 10000
Set-LiteralText $ws.Cells.Item(9, 3) "2025-11-10"

$ws.Range("A10").Value = "Groceries"
$ws.Range("B10").Value = 2000
Set-LiteralText $ws.Cells.Item(10, 3) "2025-11-08"
